$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill the "A" (sequence number) column for existing rows 20-24 ---
# (Diff shows these previously-blank A cells gaining plain numeric values.)
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23

# --- New Q&A rows 25-33 ---
# Rows 25 & 26 are "question" style rows (same look as B20/B21 -- Segoe UI 12,
# colour FF343541, row height 17.25). Rows 27-33 are plain "answer" rows.

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 'Milyen telepíthető auth-token lehetőségek vannak?'

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 'Hogyan adok hozzá Django REST framework TokenAuthentication-t a meglévő programhoz?'

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 'Django-REST-Knox van telepítve. Kérem a view.py és az url.py részeinek megírását,  az API-n keresztül  regisztrálás, login és logout működjöm.'

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 'Az előző megoldás esetén a  jelszó olvasható módon került tárolásra regisztráláskor nem lett HASH-elve, mit kell változtatni?'

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 'Írd át a LoginView és a LogoutView részeket, hogy a HASH jelszavakat kezelje:'

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 'Mit  kell megadni az előzőekben létrehozottban a sikeres logaout-hoz?'

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 'Milyen adatokat kell megadni az API POST kéréshez?'

# Row 32 has no A value in the target workbook.
$ws.Range("B32").Value = 'Alakítsd át úgy, hogy kártya "card" formában legyen a bejelentkezés form. Kékből zöldbe bal felső sarokból jobba színátmenetes háttérrel:'

# Row 33 has no A value either.
$ws.Range("B33").Value = 'Ez alőzőben létrehozott kártya az input mezők szélességéhez képest 20%-al legyen nagyobb!'

# --- Formatting: make B25/B26 look like the other "question" rows (e.g. B20/B21) ---
# Copy/paste the cell *format only* from an existing question row so the
# engine reuses the existing style entry instead of synthesising new
# (duplicate) font/style records.
$ws.Range("B20").Copy()
$ws.Range("B25").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B21").Copy()
$ws.Range("B26").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Rows(25).RowHeight = 17.25
$ws.Rows(26).RowHeight = 17.25

$ws.Range("B33").Select()
